# Daily attendance processing - 2026-02-14 16:35:16 UTC
# Swap "Miss Dina Nasr, Administrator" -> "Administrator, Miss Dina Nasr"
# in the "Recorded By" column (G) of the session analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "Miss Dina Nasr, Administrator"
$newText = "Administrator, Miss Dina Nasr"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
